$d = $word.ActiveDocument

# The original table (the one SPPR signature/table block) is deleted first so
# that $d.Paragraphs.Count afterwards reflects only the original top-level
# body paragraphs (table-cell paragraphs are otherwise included in the count).
$d.Tables(1).Delete()
$origParaCount = $d.Paragraphs.Count

$newBodyXml = @'
<w:p>
      <w:pPr>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>{{kop_surat}}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:t>______________________________________________________________________</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="24"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>SURAT PERINTAH PENDEBITAN REKENING (SPPR)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:t>Nomor: {{nomor_dokumen}}    Tanggal: {{tanggal_dokumen:tanggal}}</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Kepada Yth.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Pimpinan {{nama_bank}}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>di tempat</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t>Saya yang bertanda tangan di bawah ini selaku Pejabat Pembuat Komitmen atas nama Kuasa Pengguna Anggaran, memerintahkan agar dilakukan pendebitan rekening menggunakan kartu debit dengan keterangan sebagai berikut:</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:tbl>
      <w:tblPr>
        <w:tblStyle w:val="TableGrid"/>
        <w:tblW w:type="auto" w:w="0"/>
        <w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="4844"/>
        <w:gridCol w:w="4844"/>
      </w:tblGrid>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="2160"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Nomor Rekening</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="6480"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>: {{nomor_rekening}}</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="2160"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Nama Rekening</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="6480"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>: {{nama_rekening}}</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="2160"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Sejumlah</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="6480"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>: Rp {{nilai:rupiah}}</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="2160"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Terbilang</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="6480"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>: {{nilai:terbilang}}</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="2160"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Hari/Tanggal</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="6480"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>: {{tanggal_pencairan:tanggal}}</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Berkenaan dengan hal tersebut, mohon bantuan Saudara untuk membantu kelancaran transaksi dimaksud.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Demikian disampaikan, atas bantuan dan kerja sama yang baik diucapkan terima kasih.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:tbl>
      <w:tblPr>
        <w:tblW w:type="auto" w:w="0"/>
        <w:jc w:val="right"/>
        <w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="9689"/>
      </w:tblGrid>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="9689"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
            </w:pPr>
            <w:r>
              <w:t>{{kota}}, {{tanggal_dokumen:tanggal}}</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="9689"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
            </w:pPr>
            <w:r>
              <w:t>Pejabat Pembuat Komitmen</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="9689"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:br/>
              <w:br/>
              <w:br/>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="9689"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
              </w:rPr>
              <w:t>{{nama_ppk}}</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
            </w:pPr>
            <w:r>
              <w:t>NIP. {{nip_ppk}}</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
'@

$pkgHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$fullXml = $pkgHeader + $newBodyXml + $pkgFooter

# Insert the replacement content at the very end of the document (just
# before the final section break), then strip out everything that preceded
# it -- i.e. the original template paragraphs.
$insertPoint = $d.Content
$insertPoint.Collapse(0)
$insertPoint.InsertXML($fullXml)

for ($i = 1; $i -le $origParaCount; $i++) {
    $d.Paragraphs(1).Range.Delete()
}

# Fix the right page margin (1417 -> 1134 twips = 70.85pt -> 56.7pt)
$d.PageSetup.RightMargin = 56.7

Write-Host "Done. Paragraphs:" $d.Paragraphs.Count "Tables:" $d.Tables.Count
